# Apply the diff: replace the arithmetic problems in the table with new ones.
# Each pair is unique (no old/new collisions), so straightforward sequential
# Find & Replace across the whole document content is safe.

$d = $word.ActiveDocument

$replacements = @(
    @("798×4=", "572×3="),
    @("895×8=", "795×8="),
    @("979×9=", "542×6="),
    @("460×6=", "739×7="),
    @("191×8=", "980×4="),
    @("397×2=", "781×4="),
    @("108×2=", "177×3="),
    @("321×9=", "195×3="),
    @("981×4=", "111×2="),
    @("527×5=", "388×6="),
    @("252×6=", "578×2="),
    @("976×5=", "887×7="),
    @("133×8=", "937×9="),
    @("523×8=", "582×6="),
    @("125×7=", "136×9="),
    @("794×6=", "985×3="),
    @("133×5=", "590×4="),
    @("279×5=", "333×7="),
    @("171×2=", "108×4="),
    @("678×6=", "529×2="),
    @("457×5=", "101×7="),
    @("670×2=", "232×9="),
    @("250×6=", "814×9="),
    @("124×6=", "806×6="),
    @("474×2=", "270×2=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
